$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 225, shifting existing rows 225:322 down to 226:323.
$ws.Rows(225).Insert()

# Populate the newly inserted row 225 with the new weekly price record.
$ws.Range("A225").Value = 7
$ws.Range("B225").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C225").Value = "Ñuble"
$ws.Range("D225").Value = 44924
$ws.Range("E225").Value = 16
$ws.Range("F225").Value = 100112003
$ws.Range("G225").Value = "Ajo"
$ws.Range("H225").Value = "Chino"
$ws.Range("I225").Value = "Primera"
$ws.Range("J225").Value = 60
$ws.Range("K225").Value = 17000
$ws.Range("L225").Value = 18000
$ws.Range("M225").Value = 17500
$ws.Range("N225").Value = "`$/malla 10 kilos"
$ws.Range("O225").Value = "China"
$ws.Range("P225").Value = 1750
$ws.Range("Q225").Value = 10
$ws.Range("R225").Value = "Hortaliza"
